$wb = $excel.ActiveWorkbook

# The edit was made on the "Repayment Schedule" sheet: the user inserted a new
# (blank) column before column N, shifting the old "Late" / "Heading" /
# "Outstanding" columns one place to the right, and left the sheet active
# with cell R8 selected.
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate() | Out-Null

$ws.Columns("N").Insert() | Out-Null

# A freshly inserted column picks up the default sheet width; give it the
# same width as its left neighbour (column M), matching Excel's own
# behaviour when a column is inserted in the middle of a formatted table.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

$ws.Range("R8").Select() | Out-Null
